$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 50.4
$ws.Cells.Item(5, 9).Value = 50.4
$ws.Cells.Item(5, 11).Value = 50.4
$ws.Cells.Item(5, 13).Value = 64.59999999999999
$ws.Cells.Item(6, 8).Value = 939.8333
$ws.Cells.Item(6, 9).Value = 1007.0909
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 3021.2727
$ws.Cells.Item(6, 12).Value = 600
$ws.Cells.Item(6, 13).Value = -2909.2727
$ws.Cells.Item(6, 14).Value = -824
$ws.Cells.Item(38, 8).Value = 387.25
$ws.Cells.Item(38, 9).Value = 396.66666
$ws.Cells.Item(38, 10).Value = 359
$ws.Cells.Item(38, 11).Value = 1189.99998
$ws.Cells.Item(38, 12).Value = 1077
$ws.Cells.Item(38, 13).Value = -817.9999800000001
$ws.Cells.Item(38, 14).Value = -1821
$ws.Cells.Item(40, 8).Value = 3474.1292
$ws.Cells.Item(40, 9).Value = 1124.5
$ws.Cells.Item(40, 11).Value = 1124.5
$ws.Cells.Item(40, 13).Value = -949.5
$ws.Cells.Item(58, 8).Value = 70.75
$ws.Cells.Item(58, 9).Value = 70.75
$ws.Cells.Item(58, 11).Value = 212.25
$ws.Cells.Item(58, 13).Value = -62.25
$ws.Cells.Item(86, 8).Value = 375126600
$ws.Cells.Item(86, 9).Value = 500000700
$ws.Cells.Item(86, 10).Value = 250252500
$ws.Cells.Item(86, 11).Value = 500000700
$ws.Cells.Item(86, 12).Value = 250252500
$ws.Cells.Item(86, 13).Value = -499999577
$ws.Cells.Item(86, 14).Value = -250254746
$ws.Cells.Item(89, 8).Value = 375126600
$ws.Cells.Item(89, 9).Value = 500000700
$ws.Cells.Item(89, 10).Value = 250252500
$ws.Cells.Item(89, 11).Value = 2500003500
$ws.Cells.Item(89, 12).Value = 1251262500
$ws.Cells.Item(89, 13).Value = -2499997884
$ws.Cells.Item(89, 14).Value = -1251273732
$ws.Cells.Item(98, 8).Value = 5697.4614
$ws.Cells.Item(98, 9).Value = 5132.4443
$ws.Cells.Item(98, 11).Value = 5132.4443
$ws.Cells.Item(98, 13).Value = -3634.4443
$ws.Cells.Item(112, 8).Value = 93164.73
$ws.Cells.Item(112, 9).Value = 1495
$ws.Cells.Item(112, 11).Value = 4485
$ws.Cells.Item(112, 13).Value = -3377
$ws.Cells.Item(113, 8).Value = 4220.1816
$ws.Cells.Item(113, 9).Value = 3158
$ws.Cells.Item(113, 11).Value = 3158
$ws.Cells.Item(113, 13).Value = 96
$ws.Cells.Item(116, 8).Value = 10298.954
$ws.Cells.Item(116, 9).Value = 11727.111
$ws.Cells.Item(116, 11).Value = 11727.111
$ws.Cells.Item(116, 13).Value = -8285.111000000001
$ws.Cells.Item(122, 8).Value = 5697.4614
$ws.Cells.Item(122, 9).Value = 5132.4443
$ws.Cells.Item(122, 11).Value = 15397.3329
$ws.Cells.Item(122, 13).Value = -12947.3329
$ws.Cells.Item(129, 8).Value = 2729.75
$ws.Cells.Item(129, 9).Value = 2246
$ws.Cells.Item(129, 10).Value = 2826.5
$ws.Cells.Item(129, 11).Value = 6738
$ws.Cells.Item(129, 12).Value = 8479.5
$ws.Cells.Item(129, 13).Value = -1738
$ws.Cells.Item(129, 14).Value = -18479.5
$ws.Cells.Item(132, 8).Value = 1482.6129
$ws.Cells.Item(132, 9).Value = 1533.4828
$ws.Cells.Item(132, 11).Value = 4600.4484
$ws.Cells.Item(132, 13).Value = -2070.4484
$ws.Cells.Item(135, 8).Value = 930.86664
$ws.Cells.Item(135, 9).Value = 930.86664
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 8377.79976
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 13).Value = -5842.79976
$ws.Cells.Item(135, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 3317.0708
$ws.Cells.Item(138, 9).Value = 1991.9286
$ws.Cells.Item(138, 10).Value = 3839.6619
$ws.Cells.Item(138, 11).Value = 5975.7858
$ws.Cells.Item(138, 12).Value = 11518.9857
$ws.Cells.Item(138, 13).Value = -835.7857999999997
$ws.Cells.Item(138, 14).Value = -21798.9857
$ws.Cells.Item(141, 8).Value = 10779.8
$ws.Cells.Item(141, 9).Value = 10750
$ws.Cells.Item(141, 10).Value = 10799.667
$ws.Cells.Item(141, 11).Value = 32250
$ws.Cells.Item(141, 12).Value = 32399.001
$ws.Cells.Item(141, 13).Value = -27070
$ws.Cells.Item(141, 14).Value = -42759.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 23810598
$ws.Cells.Item(32, 9).Value = 25000876
$ws.Cells.Item(32, 11).Value = 25000876
$ws.Cells.Item(32, 13).Value = -25000589
$ws.Cells.Item(45, 8).Value = 4782.3887
$ws.Cells.Item(45, 9).Value = 4470.2144
$ws.Cells.Item(45, 10).Value = 5875
$ws.Cells.Item(45, 11).Value = 4470.2144
$ws.Cells.Item(45, 12).Value = 5875
$ws.Cells.Item(45, 13).Value = -4093.2144
$ws.Cells.Item(45, 14).Value = -6629
$ws.Cells.Item(61, 8).Value = 3358.8635
$ws.Cells.Item(61, 9).Value = 3183.111
$ws.Cells.Item(61, 10).Value = 4149.75
$ws.Cells.Item(61, 11).Value = 3183.111
$ws.Cells.Item(61, 12).Value = 4149.75
$ws.Cells.Item(61, 13).Value = -2971.111
$ws.Cells.Item(61, 14).Value = -4573.75
$ws.Cells.Item(74, 8).Value = 1904.826
$ws.Cells.Item(74, 9).Value = 1991.0476
$ws.Cells.Item(74, 11).Value = 1991.0476
$ws.Cells.Item(74, 13).Value = -1117.0476
$ws.Cells.Item(77, 8).Value = 1904.826
$ws.Cells.Item(77, 9).Value = 1991.0476
$ws.Cells.Item(77, 11).Value = 9955.238000000001
$ws.Cells.Item(77, 13).Value = -5587.238000000001
$ws.Cells.Item(97, 8).Value = 2432.0833
$ws.Cells.Item(97, 9).Value = 1407.4
$ws.Cells.Item(97, 11).Value = 1407.4
$ws.Cells.Item(97, 13).Value = -911.4000000000001
$ws.Cells.Item(110, 8).Value = 1754.0555
$ws.Cells.Item(110, 9).Value = 813.3077
$ws.Cells.Item(110, 11).Value = 813.3077
$ws.Cells.Item(110, 13).Value = 1231.6923
$ws.Cells.Item(122, 8).Value = 2227.258
$ws.Cells.Item(122, 9).Value = 1606.125
$ws.Cells.Item(122, 10).Value = 4356.857
$ws.Cells.Item(122, 11).Value = 4818.375
$ws.Cells.Item(122, 12).Value = 13070.571
$ws.Cells.Item(122, 13).Value = -2368.375
$ws.Cells.Item(122, 14).Value = -17970.571
$ws.Cells.Item(132, 8).Value = 2025.32
$ws.Cells.Item(132, 9).Value = 1957.4043
$ws.Cells.Item(132, 10).Value = 3089.3333
$ws.Cells.Item(132, 11).Value = 5872.2129
$ws.Cells.Item(132, 12).Value = 9267.999899999999
$ws.Cells.Item(132, 13).Value = -3342.2129
$ws.Cells.Item(132, 14).Value = -14327.9999
$ws.Cells.Item(133, 8).Value = 65000
$ws.Cells.Item(133, 10).Value = 65000
$ws.Cells.Item(133, 12).Value = 65000
$ws.Cells.Item(133, 14).Value = -70060
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 50000
$ws.Cells.Item(135, 10).Value = 50000
$ws.Cells.Item(135, 12).Value = 50000
$ws.Cells.Item(135, 14).Value = -60140
$ws.Cells.Item(136, 8).Value = 3358.8635
$ws.Cells.Item(136, 9).Value = 3183.111
$ws.Cells.Item(136, 10).Value = 4149.75
$ws.Cells.Item(136, 11).Value = 9549.332999999999
$ws.Cells.Item(136, 12).Value = 12449.25
$ws.Cells.Item(136, 13).Value = -6999.332999999999
$ws.Cells.Item(136, 14).Value = -17549.25
$ws.Cells.Item(138, 8).Value = 93332.664
$ws.Cells.Item(138, 10).Value = 93332.664
$ws.Cells.Item(138, 12).Value = 93332.664
$ws.Cells.Item(138, 14).Value = -103612.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 2499.5
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 2499.5
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 2499.5
$ws.Cells.Item(12, 13).ClearContents()
$ws.Cells.Item(12, 14).Value = -2835.5
$ws.Cells.Item(86, 8).Value = 3197.5
$ws.Cells.Item(86, 9).Value = 2965
$ws.Cells.Item(86, 10).Value = 3430
$ws.Cells.Item(86, 11).Value = 2965
$ws.Cells.Item(86, 12).Value = 3430
$ws.Cells.Item(86, 13).Value = -1842
$ws.Cells.Item(86, 14).Value = -5676
$ws.Cells.Item(89, 8).Value = 3197.5
$ws.Cells.Item(89, 9).Value = 2965
$ws.Cells.Item(89, 10).Value = 3430
$ws.Cells.Item(89, 11).Value = 14825
$ws.Cells.Item(89, 12).Value = 17150
$ws.Cells.Item(89, 13).Value = -9209
$ws.Cells.Item(89, 14).Value = -28382
$ws.Cells.Item(94, 8).Value = 425
$ws.Cells.Item(94, 9).Value = 487.73077
$ws.Cells.Item(94, 10).Value = 192
$ws.Cells.Item(94, 11).Value = 487.73077
$ws.Cells.Item(94, 12).Value = 192
$ws.Cells.Item(94, 13).Value = -36.73077000000001
$ws.Cells.Item(94, 14).Value = -1094
$ws.Cells.Item(99, 8).Value = 2607.0476
$ws.Cells.Item(99, 9).Value = 2171.4285
$ws.Cells.Item(99, 11).Value = 2171.4285
$ws.Cells.Item(99, 13).Value = -673.4285
$ws.Cells.Item(105, 8).Value = 5307.8335
$ws.Cells.Item(105, 9).Value = 5169.4
$ws.Cells.Item(105, 11).Value = 5169.4
$ws.Cells.Item(105, 13).Value = -3422.4
$ws.Cells.Item(134, 8).Value = 1962.7949
$ws.Cells.Item(134, 9).Value = 1486.091
$ws.Cells.Item(134, 10).Value = 4584.6665
$ws.Cells.Item(134, 11).Value = 4458.272999999999
$ws.Cells.Item(134, 12).Value = 13753.9995
$ws.Cells.Item(134, 13).Value = -1923.272999999999
$ws.Cells.Item(134, 14).Value = -18823.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 23
$ws.Cells.Item(7, 10).Value = 13
$ws.Cells.Item(7, 12).Value = 13
$ws.Cells.Item(7, 14).Value = -239
$ws.Cells.Item(22, 8).Value = 4290.0884
$ws.Cells.Item(22, 9).Value = 3923.5833
$ws.Cells.Item(22, 11).Value = 3923.5833
$ws.Cells.Item(22, 13).Value = -3573.5833
$ws.Cells.Item(31, 8).Value = 2561.6382
$ws.Cells.Item(31, 9).Value = 1445.625
$ws.Cells.Item(31, 11).Value = 1445.625
$ws.Cells.Item(31, 13).Value = -1150.625
$ws.Cells.Item(34, 8).Value = 2561.6382
$ws.Cells.Item(34, 9).Value = 1445.625
$ws.Cells.Item(34, 11).Value = 1445.625
$ws.Cells.Item(34, 13).Value = -1243.625
$ws.Cells.Item(58, 8).Value = 3707.4482
$ws.Cells.Item(58, 9).Value = 2784.2856
$ws.Cells.Item(58, 10).Value = 4001.182
$ws.Cells.Item(58, 11).Value = 2784.2856
$ws.Cells.Item(58, 12).Value = 4001.182
$ws.Cells.Item(58, 13).Value = -2581.2856
$ws.Cells.Item(58, 14).Value = -4407.182
$ws.Cells.Item(99, 8).Value = 3952.8
$ws.Cells.Item(104, 8).Value = 79980
$ws.Cells.Item(104, 10).Value = 79980
$ws.Cells.Item(104, 12).Value = 79980
$ws.Cells.Item(104, 14).Value = -85222
$ws.Cells.Item(126, 8).Value = 3952.8
$ws.Cells.Item(132, 8).Value = 6123.9
$ws.Cells.Item(132, 9).Value = 6137.6665
$ws.Cells.Item(132, 11).Value = 18412.9995
$ws.Cells.Item(132, 13).Value = -15882.9995
$ws.Cells.Item(136, 8).Value = 3707.4482
$ws.Cells.Item(136, 9).Value = 2784.2856
$ws.Cells.Item(136, 10).Value = 4001.182
$ws.Cells.Item(136, 11).Value = 8352.856800000001
$ws.Cells.Item(136, 12).Value = 12003.546
$ws.Cells.Item(136, 13).Value = -5802.856800000001
$ws.Cells.Item(136, 14).Value = -17103.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 841.8
$ws.Cells.Item(14, 9).Value = 841.8
$ws.Cells.Item(14, 11).Value = 2525.4
$ws.Cells.Item(14, 13).Value = -2352.4
$ws.Cells.Item(51, 8).Value = 1434.6666
$ws.Cells.Item(51, 9).Value = 821.6
$ws.Cells.Item(51, 11).Value = 2464.8
$ws.Cells.Item(51, 13).Value = -2004.8
$ws.Cells.Item(68, 8).Value = 1772.6923
$ws.Cells.Item(68, 9).Value = 1252.5883
$ws.Cells.Item(68, 10).Value = 2755.111
$ws.Cells.Item(68, 11).Value = 3757.7649
$ws.Cells.Item(68, 12).Value = 8265.332999999999
$ws.Cells.Item(68, 13).Value = -2946.7649
$ws.Cells.Item(68, 14).Value = -9887.332999999999
$ws.Cells.Item(71, 8).Value = 1772.6923
$ws.Cells.Item(71, 9).Value = 1252.5883
$ws.Cells.Item(71, 10).Value = 2755.111
$ws.Cells.Item(71, 11).Value = 11273.2947
$ws.Cells.Item(71, 12).Value = 24795.999
$ws.Cells.Item(71, 13).Value = -7217.294699999999
$ws.Cells.Item(71, 14).Value = -32907.999
$ws.Cells.Item(88, 8).Value = 3800
$ws.Cells.Item(88, 10).Value = 3800
$ws.Cells.Item(88, 12).Value = 11400
$ws.Cells.Item(88, 14).Value = -12256
$ws.Cells.Item(91, 8).Value = 3800
$ws.Cells.Item(91, 10).Value = 3800
$ws.Cells.Item(91, 12).Value = 11400
$ws.Cells.Item(91, 14).Value = -14364
$ws.Cells.Item(112, 8).Value = 5533
$ws.Cells.Item(112, 10).Value = 5799.5
$ws.Cells.Item(112, 12).Value = 17398.5
$ws.Cells.Item(112, 14).Value = -19614.5
$ws.Cells.Item(122, 8).Value = 491.45456
$ws.Cells.Item(122, 9).Value = 307.125
$ws.Cells.Item(122, 10).Value = 983
$ws.Cells.Item(122, 11).Value = 2764.125
$ws.Cells.Item(122, 12).Value = 8847
$ws.Cells.Item(122, 13).Value = -314.125
$ws.Cells.Item(122, 14).Value = -13747
$ws.Cells.Item(129, 8).Value = 2186.16
$ws.Cells.Item(129, 10).Value = 2495.1904
$ws.Cells.Item(129, 12).Value = 7485.5712
$ws.Cells.Item(129, 14).Value = -17485.5712

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3044.875
$ws.Cells.Item(80, 9).Value = 2825.6
$ws.Cells.Item(80, 11).Value = 2825.6
$ws.Cells.Item(80, 13).Value = -1827.6
$ws.Cells.Item(83, 8).Value = 3044.875
$ws.Cells.Item(83, 9).Value = 2825.6
$ws.Cells.Item(83, 11).Value = 14128
$ws.Cells.Item(83, 13).Value = -9136
$ws.Cells.Item(113, 8).Value = 12145.366
$ws.Cells.Item(113, 9).Value = 2170.76
$ws.Cells.Item(113, 10).Value = 27730.688
$ws.Cells.Item(113, 11).Value = 2170.76
$ws.Cells.Item(113, 12).Value = 27730.688
$ws.Cells.Item(113, 13).Value = -0.7600000000002183
$ws.Cells.Item(113, 14).Value = -32070.688
$ws.Cells.Item(122, 8).Value = 2944.5217
$ws.Cells.Item(122, 9).Value = 2174.7273
$ws.Cells.Item(122, 11).Value = 6524.1819
$ws.Cells.Item(122, 13).Value = -4074.1819
$ws.Cells.Item(126, 8).Value = 3215.6086
$ws.Cells.Item(126, 9).Value = 2597.3076
$ws.Cells.Item(126, 11).Value = 7791.9228
$ws.Cells.Item(126, 13).Value = -5321.9228
$ws.Cells.Item(132, 8).Value = 1034.1666
$ws.Cells.Item(132, 9).Value = 777.88464
$ws.Cells.Item(132, 10).Value = 2700
$ws.Cells.Item(132, 11).Value = 2333.65392
$ws.Cells.Item(132, 12).Value = 8100
$ws.Cells.Item(132, 13).Value = 196.3460800000003
$ws.Cells.Item(132, 14).Value = -13160

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3531.3333
$ws.Cells.Item(40, 9).Value = 3454.4
$ws.Cells.Item(40, 10).Value = 3627.5
$ws.Cells.Item(40, 11).Value = 3454.4
$ws.Cells.Item(40, 12).Value = 3627.5
$ws.Cells.Item(40, 13).Value = -3318.4
$ws.Cells.Item(40, 14).Value = -3899.5
$ws.Cells.Item(46, 8).Value = 6359
$ws.Cells.Item(46, 10).Value = 6845.241
$ws.Cells.Item(46, 12).Value = 6845.241
$ws.Cells.Item(46, 14).Value = -7221.241
$ws.Cells.Item(81, 8).Value = 114000
$ws.Cells.Item(81, 10).Value = 114000
$ws.Cells.Item(81, 12).Value = 114000
$ws.Cells.Item(81, 14).Value = -115996
$ws.Cells.Item(82, 8).Value = 3299.0625
$ws.Cells.Item(82, 9).Value = 4001
$ws.Cells.Item(82, 11).Value = 4001
$ws.Cells.Item(82, 13).Value = -3640
$ws.Cells.Item(84, 8).Value = 114000
$ws.Cells.Item(84, 10).Value = 114000
$ws.Cells.Item(84, 12).Value = 342000
$ws.Cells.Item(84, 14).Value = -351984
$ws.Cells.Item(85, 8).Value = 3299.0625
$ws.Cells.Item(85, 9).Value = 4001
$ws.Cells.Item(85, 11).Value = 4001
$ws.Cells.Item(85, 13).Value = -2753
$ws.Cells.Item(93, 8).Value = 55557284
$ws.Cells.Item(93, 9).Value = 83334820
$ws.Cells.Item(93, 11).Value = 83334820
$ws.Cells.Item(93, 13).Value = -83333572
$ws.Cells.Item(132, 8).Value = 6298.769
$ws.Cells.Item(132, 9).Value = 6859.522
$ws.Cells.Item(132, 10).Value = 1999.6666
$ws.Cells.Item(132, 11).Value = 20578.566
$ws.Cells.Item(132, 12).Value = 5998.9998
$ws.Cells.Item(132, 13).Value = -18048.566
$ws.Cells.Item(132, 14).Value = -11058.9998
$ws.Cells.Item(136, 8).Value = 28679.785
$ws.Cells.Item(136, 9).Value = 28679.785
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 86039.355
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -83489.355
$ws.Cells.Item(136, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 7100.3335
$ws.Cells.Item(62, 9).Value = 7438
$ws.Cells.Item(62, 11).Value = 7438
$ws.Cells.Item(62, 13).Value = -6814
$ws.Cells.Item(65, 8).Value = 7100.3335
$ws.Cells.Item(65, 9).Value = 7438
$ws.Cells.Item(65, 11).Value = 37190
$ws.Cells.Item(65, 13).Value = -34070
$ws.Cells.Item(70, 8).Value = 29071.428
$ws.Cells.Item(70, 9).Value = 17500
$ws.Cells.Item(70, 11).Value = 17500
$ws.Cells.Item(70, 13).Value = -17185
$ws.Cells.Item(73, 8).Value = 29071.428
$ws.Cells.Item(73, 9).Value = 17500
$ws.Cells.Item(73, 11).Value = 17500
$ws.Cells.Item(73, 13).Value = -16408
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 961
$ws.Cells.Item(113, 9).Value = 912.25
$ws.Cells.Item(113, 11).Value = 2736.75
$ws.Cells.Item(113, 13).Value = -566.75
$ws.Cells.Item(122, 8).Value = 3907.7568
$ws.Cells.Item(122, 9).Value = 3673.6177
$ws.Cells.Item(122, 10).Value = 6561.3335
$ws.Cells.Item(122, 11).Value = 11020.8531
$ws.Cells.Item(122, 12).Value = 19684.0005
$ws.Cells.Item(122, 13).Value = -8570.8531
$ws.Cells.Item(122, 14).Value = -24584.0005
$ws.Cells.Item(132, 8).Value = 3239.5386
$ws.Cells.Item(132, 9).Value = 3374
$ws.Cells.Item(132, 11).Value = 10122
$ws.Cells.Item(132, 13).Value = -7592
$ws.Cells.Item(136, 8).Value = 28903.078
$ws.Cells.Item(136, 9).Value = 2221.5757
$ws.Cells.Item(136, 11).Value = 6664.7271
$ws.Cells.Item(136, 13).Value = -4114.7271

